# Atualização de bases das ligas, do dia: 20-02-2024 às 23:00
#
# The underlying match-id (column B) bookkeeping got re-synced with the
# upstream feed, which re-ordered a handful of fixture rows. Column A
# (the running row index) and columns C/D/E (Div / Div Original Name /
# Date) stay put; every other column (B, F..AC) for the affected rows
# is swapped/rotated to its corrected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple pairwise swaps (rows 27/28, 101/102, 164/165) ---

$v27 = $ws.Range("B27:AC27").Value2
$v28 = $ws.Range("B28:AC28").Value2
$ws.Range("B27:AC27").Value2 = $v28
$ws.Range("B28:AC28").Value2 = $v27

$v101 = $ws.Range("B101:AC101").Value2
$v102 = $ws.Range("B102:AC102").Value2
$ws.Range("B101:AC101").Value2 = $v102
$ws.Range("B102:AC102").Value2 = $v101

$v164 = $ws.Range("B164:AC164").Value2
$v165 = $ws.Range("B165:AC165").Value2
$ws.Range("B164:AC164").Value2 = $v165
$ws.Range("B165:AC165").Value2 = $v164

# --- Cyclic rotation across rows 175-179 ---
#   data at row 175 -> row 179
#   data at row 176 -> row 177
#   data at row 177 -> row 178
#   data at row 178 -> row 175
#   data at row 179 -> row 176

$v175 = $ws.Range("B175:AC175").Value2
$v176 = $ws.Range("B176:AC176").Value2
$v177 = $ws.Range("B177:AC177").Value2
$v178 = $ws.Range("B178:AC178").Value2
$v179 = $ws.Range("B179:AC179").Value2

$ws.Range("B179:AC179").Value2 = $v175
$ws.Range("B177:AC177").Value2 = $v176
$ws.Range("B178:AC178").Value2 = $v177
$ws.Range("B175:AC175").Value2 = $v178
$ws.Range("B176:AC176").Value2 = $v179
